$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to Text format so numeric-looking strings
# (e.g. "576.23") are stored as text, matching the source data which
# keeps all Price values as plain strings (not numbers).
foreach ($ref in @("D5","D6","D8","D10","D19","D20","D21","D24","D25","D27","D29","D31","D32","D35","D36","D37","D38","D39","D40","D41","D42","D43","D46","D47","D48","D49","D50")) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "67.184.89"
$ws.Range("E2").Value = "  +3.72%  "

# Row 3
$ws.Range("D3").Value = "3.229.36"

# Row 5
$ws.Range("D5").Value = "576.23"
$ws.Range("E5").Value = "  +1.90%  "

# Row 6
$ws.Range("D6").Value = "180.37"
$ws.Range("E6").Value = "  +5.42%  "

# Row 8
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -4.40%  "

# Row 9
$ws.Range("D9").Value = "3.226.75"
$ws.Range("E9").Value = "  +1.86%  "

# Row 10
$ws.Range("D10").Value = "0.129"
$ws.Range("E10").Value = "  +3.18%  "

# Row 11
$ws.Range("E11").Value = "  +3.15%  "

# Row 12
$ws.Range("E12").Value = "  +4.33%  "

# Row 13
$ws.Range("D13").Value = "3.786.99"
$ws.Range("E13").Value = "  +1.70%  "

# Row 14
$ws.Range("E14").Value = "  +1.26%  "

# Row 15
$ws.Range("E15").Value = "  +1.52%  "

# Row 16
$ws.Range("D16").Value = "67.126.30"
$ws.Range("E16").Value = "  +3.86%  "

# Row 17
$ws.Range("E17").Value = "  +2.20%  "

# Row 18
$ws.Range("D18").Value = "3.234.35"
$ws.Range("E18").Value = "  +2.00%  "

# Row 19
$ws.Range("D19").Value = "5.78"
$ws.Range("E19").Value = "  +0.94%  "

# Row 20
$ws.Range("D20").Value = "13.37"
$ws.Range("E20").Value = "  +2.87%  "

# Row 21
$ws.Range("D21").Value = "373.27"
$ws.Range("E21").Value = "  +4.80%  "

# Row 22
$ws.Range("E22").Value = "  +4.05%  "

# Row 23
$ws.Range("E23").Value = "  -0.71%  "

# Row 24
$ws.Range("D24").Value = "70.91"
$ws.Range("E24").Value = "  +3.33%  "

# Row 25
$ws.Range("D25").Value = "0.508"
$ws.Range("E25").Value = "  +1.29%  "

# Row 26
$ws.Range("E26").Value = "  +0.48%  "

# Row 27
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  -0.91%  "

# Row 28
$ws.Range("E28").Value = "  +2.82%  "

# Row 29
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.36%  "

# Row 30
$ws.Range("E30").Value = "  +3.70%  "

# Row 31
$ws.Range("D31").Value = "5.63"
$ws.Range("E31").Value = "  +3.54%  "

# Row 32
$ws.Range("D32").Value = "22.50"
$ws.Range("E32").Value = "  +2.12%  "

# Row 33
$ws.Range("E33").Value = "  +0.02%  "

# Row 34
$ws.Range("E34").Value = "  +3.71%  "

# Row 35
$ws.Range("D35").Value = "6.83"
$ws.Range("E35").Value = "  +2.13%  "

# Row 36
$ws.Range("D36").Value = "162.41"
$ws.Range("E36").Value = "  +5.81%  "

# Row 37
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  +2.91%  "

# Row 38
$ws.Range("D38").Value = "0.852"
$ws.Range("E38").Value = "  +2.81%  "

# Row 39
$ws.Range("D39").Value = "1.85"
$ws.Range("E39").Value = "  +6.42%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "6.75"
$ws.Range("E40").Value = "  +11.83%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.54"
$ws.Range("E41").Value = "  +0.49%  "

# Row 42
$ws.Range("D42").Value = "2.60"
$ws.Range("E42").Value = "  +3.07%  "

# Row 43
$ws.Range("D43").Value = "362.11"
$ws.Range("E43").Value = "  +12.27%  "

# Row 44
$ws.Range("E44").Value = "  +5.05%  "

# Row 45
$ws.Range("D45").Value = "2.694.64"
$ws.Range("E45").Value = "  +1.52%  "

# Row 46
$ws.Range("D46").Value = "25.48"
$ws.Range("E46").Value = "  +5.14%  "

# Row 47
$ws.Range("D47").Value = "40.54"
$ws.Range("E47").Value = "  +3.02%  "

# Row 48
$ws.Range("D48").Value = "0.0670"
$ws.Range("E48").Value = "  +2.39%  "

# Row 49
$ws.Range("D49").Value = "0.0277"
$ws.Range("E49").Value = "  +1.35%  "

# Row 50
$ws.Range("D50").Value = "0.991"
$ws.Range("E50").Value = "  +5.64%  "

# Row 51
$ws.Range("E51").Value = "  -0.31%  "
